$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute(
        $old, $false, $false, $false, $false, $false, $true, 1, $false,
        $new, 2)
    if (-not $found) {
        Write-Output "WARNING: text not found -> $old"
    }
}

# --- Change 1: "2+ years of" -> "3+ years of" ---------------------------
Replace-Text " 2+ years of" " 3+ years of"

# --- Change 2: "prioritise task" -> "prioritise task." -------------------
Replace-Text "Good time management skills and prioritise task" "Good time management skills and prioritise task."

# --- Change 3: Education section rewrite ---------------------------------
# Process bottom-to-top so that identical "old" text values (shared
# between the two education entries before the edit) never collide with
# text already written by an earlier replacement in this same pass.

# Entry 2 location line
Replace-Text "Victoria University, Sydney, NSW, Australia" "Kingston Institute Australia, Sydney, NSW, Australia"

# Entry 2 degree line
Replace-Text "Bachelor's degree in information technology: July 2020 - September 2022" "Diploma and Advanced Diploma in information technology: May 2017 - May 2020"

# Entry 1 location line
Replace-Text "QIBA, Sydney, NSW, Australia" "Victoria University, Sydney, NSW, Australia"

# Entry 1 degree line
Replace-Text "ACS Professional Year Program: March 2023 - November 2023" "Bachelor's degree in information technology: July 2020 - September 2022"
